# Daily attendance processing - 2025-12-11 15:57:30
# Swap the order of the two comma-separated names in the "Recorded By"
# column (column G) wherever "dnasr281@gmail.com" currently appears first,
# moving it to the end of the list instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -like "dnasr281@gmail.com,*") {
        $parts = $val -split ", "
        if ($parts.Count -eq 2) {
            $cell.Value2 = "$($parts[1]), $($parts[0])"
        }
    }
}
